# Weekly update: prepend two new "Camote" (Zapallo) price rows for the
# latest survey date (serial 44509), pushing the existing data block
# (rows 399:415) down by two rows (to 401:417).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first data row of the block
# (row 399). Excel shifts rows 399:415 down to 401:417, carrying their
# contents and formatting with them.
$ws.Rows("399:400").Insert()

# --- New row 399: "1a nueva(o)" ---
$ws.Cells.Item(399, 1).Value2 = 8
$ws.Cells.Item(399, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(399, 3).Value2 = "Coquimbo"
$ws.Cells.Item(399, 4).Value2 = 44509
$ws.Cells.Item(399, 5).Value2 = 4
$ws.Cells.Item(399, 6).Value2 = 100112045
$ws.Cells.Item(399, 7).Value2 = "Zapallo"
$ws.Cells.Item(399, 8).Value2 = "Camote"
$ws.Cells.Item(399, 9).Value2 = "1a nueva(o)"
$ws.Cells.Item(399, 10).Value2 = 740
$ws.Cells.Item(399, 11).Value2 = 700
$ws.Cells.Item(399, 12).Value2 = 750
$ws.Cells.Item(399, 13).Value2 = 725
$ws.Cells.Item(399, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(399, 15).Value2 = "Perú"
$ws.Cells.Item(399, 16).Value2 = 725
$ws.Cells.Item(399, 17).Value2 = 1
$ws.Cells.Item(399, 18).Value2 = "Hortaliza"

# --- New row 400: "2a nueva(o)" ---
$ws.Cells.Item(400, 1).Value2 = 8
$ws.Cells.Item(400, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(400, 3).Value2 = "Coquimbo"
$ws.Cells.Item(400, 4).Value2 = 44509
$ws.Cells.Item(400, 5).Value2 = 4
$ws.Cells.Item(400, 6).Value2 = 100112045
$ws.Cells.Item(400, 7).Value2 = "Zapallo"
$ws.Cells.Item(400, 8).Value2 = "Camote"
$ws.Cells.Item(400, 9).Value2 = "2a nueva(o)"
$ws.Cells.Item(400, 10).Value2 = 520
$ws.Cells.Item(400, 11).Value2 = 600
$ws.Cells.Item(400, 12).Value2 = 650
$ws.Cells.Item(400, 13).Value2 = 625
$ws.Cells.Item(400, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(400, 15).Value2 = "Perú"
$ws.Cells.Item(400, 16).Value2 = 625
$ws.Cells.Item(400, 17).Value2 = 1
$ws.Cells.Item(400, 18).Value2 = "Hortaliza"

# Give the new date cells the same date style ("s=2") already used by
# column D throughout the table (Insert() should carry this already,
# but make it explicit/robust).
$ws.Range("D399:D400").NumberFormat = $ws.Range("D401").NumberFormat
